$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36, shifting existing rows 36-73 down to 37-74
$ws.Rows(36).Insert()

# Populate the new row 36 with the new data record
$ws.Cells.Item(36, 1).Value = 6
$ws.Cells.Item(36, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(36, 3).Value = "Metropolitana"
$ws.Cells.Item(36, 4).Value = 45079
$ws.Cells.Item(36, 5).Value = 13
$ws.Cells.Item(36, 6).Value = 100112035
$ws.Cells.Item(36, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(36, 8).Value = "Sin especificar"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 270
$ws.Cells.Item(36, 11).Value = 18000
$ws.Cells.Item(36, 12).Value = 20000
$ws.Cells.Item(36, 13).Value = 19037
$ws.Cells.Item(36, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(36, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(36, 16).Value = 1269
$ws.Cells.Item(36, 17).Value = 15
$ws.Cells.Item(36, 18).Value = "Hortaliza"

Write-Host "done"
